# RMA End to End Test Cases Other maintenance
# Refresh the "Routing Master" lookup row with the latest engineering-item
# test record (Provar writes a new Item Number / Id pair after each run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routing Master")

$ws.Range("B2").Value = "Pro-PEItem-WLQMI"
$ws.Range("D2").Value = "a2S1K000002TLVOUA4"
